# Clean up the data dictionary:
#  - "date" is actually an ID-type field, not a safety_measures field
#  - strip stray trailing (non-breaking) spaces from the "cli", "cri", "cti"
#    and "ln_fi" variable-name entries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C9").Value = "ID"

$ws.Range("A16").Value = "cli"
$ws.Range("A17").Value = "cri"
$ws.Range("A18").Value = "cti"
$ws.Range("A21").Value = "ln_fi"

# Restore the view: scrolled so row 10 is at the top, with A22 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A22").Select()
